$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D writes to stay as text (many values look numeric, e.g. "618.67")
# by pre-formatting as Text, then clearing the format afterwards so no
# residual number-format/style is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

# --- Coin name / link swaps (Litecoin <-> Fetch.AI, Mantle <-> Filecoin) ---
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("B25").Value = "Fetch.AI"
$ws.Range("C25").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "70.164.21"
$ws.Range("D3").Value = "3.784.42"
$ws.Range("D5").Value = "618.67"
$ws.Range("D6").Value = "177.62"
$ws.Range("D7").Value = "3.778.20"
$ws.Range("D9").Value = "0.535"
$ws.Range("D10").Value = "0.172"
$ws.Range("D11").Value = "6.31"
$ws.Range("D12").Value = "0.492"
$ws.Range("D13").Value = "40.99"
$ws.Range("D14").Value = "0.0000262"
$ws.Range("D15").Value = "4.417.54"
$ws.Range("D16").Value = "3.781.78"
$ws.Range("D17").Value = "70.216.42"
$ws.Range("D19").Value = "7.61"
$ws.Range("D20").Value = "16.85"
$ws.Range("D21").Value = "511.16"
$ws.Range("D22").Value = "9.50"
$ws.Range("D23").Value = "0.727"
$ws.Range("D24").Value = "87.82"
$ws.Range("D25").Value = "2.50"
$ws.Range("D26").Value = "13.13"
$ws.Range("D27").Value = "10.99"
$ws.Range("D29").Value = "0.999"
$ws.Range("D30").Value = "2.48"
$ws.Range("D31").Value = "2.86"
$ws.Range("D33").Value = "31.40"
$ws.Range("D34").Value = "0.115"
$ws.Range("D36").Value = "1.06"
$ws.Range("D37").Value = "6.22"
$ws.Range("D39").Value = "0.333"
$ws.Range("D40").Value = "2.13"
$ws.Range("D41").Value = "51.01"
$ws.Range("D42").Value = "44.93"
$ws.Range("D43").Value = "8.74"
$ws.Range("D44").Value = "417.92"
$ws.Range("D45").Value = "2.83"
$ws.Range("D46").Value = "3.033.09"
$ws.Range("D47").Value = "0.0364"
$ws.Range("D48").Value = "27.39"
$ws.Range("D49").Value = "139.11"

# Remove the temporary Text number format so the cells end up with no
# explicit style, matching their original (unstyled) state.
$ws.Range("D2:D51").ClearFormats()

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +3.51%  "
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +5.18%  "
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("E13").Value = "  +2.80%  "
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("E15").Value = "  +3.26%  "
$ws.Range("E16").Value = "  +3.26%  "
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("E18").Value = "  -0.08%  "
$ws.Range("E19").Value = "  +1.36%  "
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("E22").Value = "  +3.21%  "
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +3.49%  "
$ws.Range("E26").Value = "  -2.70%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  +27.33%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  -1.81%  "
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("E36").Value = "  +5.41%  "
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("E38").Value = "  +5.92%  "
$ws.Range("E39").Value = "  -2.85%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("E44").Value = "  +4.84%  "
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("E46").Value = "  -4.96%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("E48").Value = "  -2.75%  "
$ws.Range("E49").Value = "  +2.37%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  +1.41%  "
